# Applies the cryptos.xlsx data refresh described by the commit:
#   "Updated cryptos list on Thu Mar  2 16:45:11 UTC 2023 with GitHub Actions"
#
# The sheet stores every data cell (Coin, Link, Price, Volume(1h)) as plain
# text, even when the text happens to look like a number (e.g. "0.9999",
# "1.001"). Writing such a string straight into a General-formatted cell
# would make Excel re-interpret it as a numeric value, which both loses the
# original textual formatting and does not match the source data. So for
# any new Price (column D) value that parses as a number we first flip the
# cell's NumberFormat to Text ("@") and then assign the value, which keeps
# it a plain text cell exactly like the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [pscustomobject]@{ Cell = 'D2'; Value = '23.315.65'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E2'; Value = '  -1.61%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D3'; Value = '1.626.47'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E3'; Value = '  -1.76%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D4'; Value = '1.003'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E4'; Value = '  +0.12%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D5'; Value = '1.002'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E5'; Value = '  +0.17%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D6'; Value = '297.87'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E6'; Value = '  -1.65%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D7'; Value = '0.3751'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E7'; Value = '  -1.62%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D8'; Value = '50.26'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E8'; Value = '  -2.44%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D9'; Value = '0.3469'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E9'; Value = '  -3.99%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D10'; Value = '0.08009'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E10'; Value = '  -2.30%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D11'; Value = '1.196'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E11'; Value = '  -2.85%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D12'; Value = '1.003'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E12'; Value = '  +0.14%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D13'; Value = '21.78'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E13'; Value = '  -3.54%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D14'; Value = '6.279'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E14'; Value = '  -3.44%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D15'; Value = '7.194'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E15'; Value = '  -2.90%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D16'; Value = '0.00001179'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E16'; Value = '  -4.18%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D17'; Value = '1.629.05'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E17'; Value = '  -1.60%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D18'; Value = '94.47'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E18'; Value = '  -3.31%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D19'; Value = '0.06942'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E19'; Value = '  -1.01%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D20'; Value = '6.589'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E20'; Value = '  -3.55%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D21'; Value = '17.22'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E21'; Value = '  -2.41%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D22'; Value = '1.002'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E22'; Value = '  +0.07%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D23'; Value = '12.30'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E23'; Value = '  -3.86%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D24'; Value = '23.333.32'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E24'; Value = '  -1.50%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D25'; Value = '2.434'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E25'; Value = '  -3.30%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D26'; Value = '3.033'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E26'; Value = '  -0.02%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D27'; Value = '20.73'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E27'; Value = '  -2.61%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D28'; Value = '151.19'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E28'; Value = '  -1.48%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D29'; Value = '5.152'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E29'; Value = '  -1.33%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D30'; Value = '131.29'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E30'; Value = '  -2.32%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D31'; Value = '1.814.28'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E31'; Value = '  -1.06%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D32'; Value = '6.675'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E32'; Value = '  -4.87%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D33'; Value = '2.131'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E33'; Value = '  -3.85%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'E34'; Value = '  -5.50%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D35'; Value = '0.9703'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E35'; Value = '  -8.18%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D36'; Value = '0.02652'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E36'; Value = '  -5.20%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D37'; Value = '0.08709'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E37'; Value = '  -0.80%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D38'; Value = '0.2407'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E38'; Value = '  -4.33%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D39'; Value = '5.794'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E39'; Value = '  -4.76%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B40'; Value = 'Hedera'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C40'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D40'; Value = '0.06640'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E40'; Value = '  -5.57%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B41'; Value = 'Aptos'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C41'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D41'; Value = '12.60'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E41'; Value = '  -3.13%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D42'; Value = '0.6761'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E42'; Value = '  -3.37%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D43'; Value = '1.284'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E43'; Value = '  -3.86%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D44'; Value = '15.34'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E44'; Value = '  -4.25%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D45'; Value = '1.001'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E45'; Value = '  +0.06%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D46'; Value = '0.6273'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E46'; Value = '  -3.71%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B47'; Value = 'PancakeSwap'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D47'; Value = '3.886'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E47'; Value = '  -1.92%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B48'; Value = 'NEARProtocol'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C48'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D48'; Value = '2.222'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E48'; Value = '  -3.88%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D49'; Value = '126.37'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E49'; Value = '  -1.23%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D50'; Value = '0.07625'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E50'; Value = '  -3.68%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D51'; Value = '1.207'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E51'; Value = '  +0.68%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
